# Goal: split the paragraph containing the M2Doc field text
#   {m:'An empty table'.emptyTable()}
# (currently stored as two runs: "{m" and ":'An empty table'.emptyTable()}")
# into four separate runs:
#   "{"
#   "m"
#   ":'An empty table'.emptyTable()"
#   "}"
#
# NOTE on this runtime's behaviour: calling Range.Text = ... or Range.Delete()
# on a sub-range that overlaps existing run text causes Word to recompute/merge
# all runs of that paragraph into the minimal run set on save. The only
# operations that reliably keep freshly-produced text as distinct <w:r> runs
# are paragraph-boundary inserts (InsertParagraphAfter/Before at a position,
# plus deleting a paragraph mark to rejoin two paragraphs). So the approach
# here is: temporarily split the target paragraph into one paragraph per
# desired run (using InsertParagraphAfter, which never touches existing run
# text), then rejoin those paragraphs by deleting the paragraph marks between
# them (which preserves run separation instead of re-merging the text).

$d = $word.ActiveDocument

# Locate the paragraph that holds the "{m: ... emptyTable() }" field text.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "{m*emptyTable*}*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the paragraph containing the '{m ... emptyTable() }' field text."
}

$p = $d.Paragraphs.Item($targetIndex)
$pStart = $p.Range.Start

# Range.Text of a paragraph includes the trailing paragraph-mark character,
# so the visible text length is one less than the string length.
$fullParaText = $p.Range.Text
$textLen = $fullParaText.Length - 1
$visibleText = $fullParaText.Substring(0, $textLen)

if (-not ($visibleText.StartsWith("{m") -and $visibleText.EndsWith("}"))) {
    throw "Unexpected paragraph content, aborting: [$visibleText]"
}

# Desired run boundaries (character offsets from the start of the paragraph):
#   [0,1)      -> "{"
#   [1,2)      -> "m"
#   [2,N-1)    -> ":'An empty table'.emptyTable()"
#   [N-1,N)    -> "}"
$splitOffsets = @(1, 2, $textLen - 1)

# Insert the paragraph breaks starting from the right-most offset first, so
# that the offsets computed above (relative to the original paragraph start)
# remain valid as we go.
$sortedDesc = $splitOffsets | Sort-Object -Descending
foreach ($off in $sortedDesc) {
    $abs = $pStart + $off
    $splitRange = $d.Range($abs, $abs)
    $splitRange.InsertParagraphAfter()
}

# We now have 4 consecutive small paragraphs (at $targetIndex .. $targetIndex+3)
# each holding one of the desired pieces as its own run. Rejoin them into a
# single paragraph by deleting the paragraph mark that separates the
# paragraph at $targetIndex from the next one, three times.
for ($j = 0; $j -lt ($splitOffsets.Count); $j++) {
    $joinPara = $d.Paragraphs.Item($targetIndex)
    $joinEnd = $joinPara.Range.End
    $markRange = $d.Range($joinEnd - 1, $joinEnd)
    $markRange.Delete()
}

$finalText = $d.Paragraphs.Item($targetIndex).Range.Text
if ($finalText -ne $fullParaText) {
    throw "Post-condition failed: paragraph text changed unexpectedly. Expected [$fullParaText] got [$finalText]"
}

"Split complete. Final paragraph text: [" + $finalText + "]"
